# The presentation's single Design ("Integral") is recoloured to the
# default Office Theme palette. In the underlying OOXML this is the
# slide master's theme part (ppt/theme/theme1.xml): its <a:clrScheme>
# is updated from the Integral colours to the standard Office colours
# (font scheme / format scheme are already identical between the two
# themes bundled with this deck, so only the 12 theme colours change).

$p   = $ppt.ActivePresentation
$sm  = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

function HexToBgrInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order matches MsoThemeColorSchemeIndex (1-12):
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000", # Text/Background - Dark 1
    "FFFFFF", # Text/Background - Light 1
    "44546A", # Text/Background - Dark 2
    "E7E6E6", # Text/Background - Light 2
    "5B9BD5", # Accent 1
    "ED7D31", # Accent 2
    "A5A5A5", # Accent 3
    "FFC000", # Accent 4
    "4472C4", # Accent 5
    "70AD47", # Accent 6
    "0563C1", # Hyperlink
    "954F72"  # Followed Hyperlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = HexToBgrInt($officeThemeColors[$i - 1])
}
